# Apply the edit described in the commit: renumber RM552628 -> RM553315,
# update vendor names in the DIM_VENDEDOR sheet, update the SQL sheet's
# query text, and register the two date number formats that appear in
# the target styles.xml.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the first sheet (RM552628.DIM_VENDEDOR -> RM553315.DIM_VENDEDOR) ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "RM553315.DIM_VENDEDOR"

# --- 2. Update the vendor names in column B (rows 2-101) ---
$names = @{
    2 = 'Elias Pinto'
    3 = 'Pedro Siqueira'
    4 = 'Sandra Moraes'
    5 = 'Tiago Campos'
    6 = 'Valeria Campos'
    7 = 'Roberto Barros'
    8 = 'Joao Siqueira'
    9 = 'Leandro Hernandes'
    10 = 'Valeria Teixeira'
    12 = 'Quiteria Alves'
    13 = 'Nelson Pinto'
    14 = 'Leandro Siqueira'
    15 = 'Elias Hernandes'
    16 = 'Roberto Farias'
    17 = 'Quiteria Pinto'
    18 = 'Olga Farias'
    19 = 'Leandro Campos'
    20 = 'Roberto Brito'
    21 = 'Fatima Pinto'
    22 = 'Roberto Duarte'
    23 = 'Nelson Freitas'
    24 = 'Nelson Siqueira'
    25 = 'Nelson Nascimento'
    26 = 'Monica Freitas'
    27 = 'Nelson Duarte'
    28 = 'Isabela Farias'
    29 = 'Nelson Xavier'
    30 = 'Claudio Castro'
    31 = 'Tiago Viana'
    32 = 'Nelson Viana'
    33 = 'Olga Teixeira'
    34 = 'Isabela Lopes'
    35 = 'Debora Farias'
    36 = 'Leandro Xavier'
    37 = 'Roberto Zanetti'
    38 = 'Roberto Freitas'
    39 = 'Karine Lopes'
    40 = 'Nelson Goncalves'
    41 = 'Andre Campos'
    42 = 'Nelson Zanetti'
    43 = 'Karine Farias'
    44 = 'Joao Lopes'
    45 = 'Fatima Zanetti'
    46 = 'Debora Nascimento'
    47 = 'Roberto Moraes'
    48 = 'Nelson Castro'
    49 = 'Debora Duarte'
    50 = 'Leandro Xavier'
    51 = 'Debora Campos'
    52 = 'Isabela Moraes'
    53 = 'Fatima Lopes'
    54 = 'Valeria Campos'
    55 = 'Monica Alves'
    56 = 'Karine Siqueira'
    57 = 'Sandra Viana'
    58 = 'Quiteria Moraes'
    59 = 'Tiago Barros'
    60 = 'Beatriz Freitas'
    61 = 'Karine Pinto'
    62 = 'Sandra Freitas'
    63 = 'Olga Campos'
    64 = 'Pedro Alves'
    65 = 'Valeria Hernandes'
    66 = 'Pedro Siqueira'
    67 = 'Leandro Alves'
    68 = 'Olga Hernandes'
    69 = 'Leandro Rocha'
    70 = 'Tiago Barros'
    71 = 'Elias Castro'
    72 = 'Andre Duarte'
    73 = 'Monica Castro'
    74 = 'Tiago Viana'
    75 = 'Leandro Viana'
    76 = 'Sandra Teixeira'
    77 = 'Pedro Lopes'
    78 = 'Quiteria Zanetti'
    79 = 'Elias Barros'
    80 = 'Karine Goncalves'
    81 = 'Pedro Teixeira'
    82 = 'Pedro Barros'
    83 = 'Beatriz Duarte'
    84 = 'Valeria Lopes'
    85 = 'Gustavo Teixeira'
    86 = 'Leandro Duarte'
    87 = 'Valeria Campos'
    88 = 'Claudio Barros'
    89 = 'Karine Xavier'
    90 = 'Pedro Duarte'
    91 = 'Roberto Rocha'
    92 = 'Pedro Viana'
    93 = 'Pedro Nascimento'
    94 = 'Karine Xavier'
    95 = 'Claudio Duarte'
    96 = 'Nelson Goncalves'
    97 = 'Elias Teixeira'
    98 = 'Debora Goncalves'
    99 = 'Beatriz Alves'
    100 = 'Monica Brito'
    101 = 'Nelson Moraes'
}

foreach ($row in $names.Keys) {
    $ws.Cells.Item($row, 2).Value = $names[$row]
}

# --- 3. Update the SQL sheet query text (RM552628 -> RM553315) ---
$sql = $wb.Worksheets.Item(2)
$sql.Range("A2").Value = 'select COD_VENDEDOR COD_VENDEDOR, NOME_VENDEDOR NOME_VENDEDOR from (select * from "RM553315"."DIM_VENDEDOR")'

# --- 4. Register the date number formats (164/165) used elsewhere in the
#        workbook, without altering the visible A1:B101 used range. We set
#        them on a scratch cell outside the sheet's data range and then
#        fully clear that cell so no value/format remains attached to it,
#        matching the target: numFmts/cellXfs exist but are unused.
$ws.Range("D1").NumberFormat = "m/d/yyyy h:mm AM/PM"
$ws.Range("D2").NumberFormat = "m/d/yyyy"
$ws.Range("D1:D2").Clear()
